$d = $word.ActiveDocument

# 1. Remove "batteries " from the paper title (keep single space before "based")
$d.Content.Find.Execute("electric vehicle (EV) batteries based on trip data", $true, $false, $false, $false, $false, $true, 1, $false, "electric vehicle (EV) based on trip data", 2)

# 2. Fix citation number: dataset reference should be [12] not [1]
$d.Content.Find.Execute("dual-Electric Vehicle Dataset (d-EVD) [1], an open-access", $true, $false, $false, $false, $false, $true, 1, $false, "dual-Electric Vehicle Dataset (d-EVD) [12], an open-access", 2)
